# Add a new installment_payment WhatsApp template (message_id 5197) as row 4
# on the "WhatsApp Templates API" sheet, matching the existing table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WhatsApp Templates API")

$msgContent = "HEADER: CHORDS MUSIC ACADEMY | BODY: Dear {Var1},`nINSTALLMENT REMINDER`n" + `
    [char]0x2022 + " Amount: " + [char]0x20B9 + "{Var2}`n" + `
    [char]0x2022 + " Due Date: {Var3}`n" + `
    [char]0x2022 + " Package: {Var4}`n" + `
    [char]0x2022 + " Expires: {Var5}`n" + `
    "Pay now to continue your musical journey.`n" + `
    [char]0x1F4DE + " 7981585309 | BUTTON: Call Us (PHONE_NUMBER) - +917981585309"

$apiUrl = "https://www.fast2sms.com/dev/whatsapp?authorization=<YOUR_API_KEY>&message_id=5197&numbers=<MOBILE_NUMBER>&variables_values=Var1|Var2|Var3|Var4|Var5"

$ws.Range("A4").Value = 5197
$ws.Range("B4").Value = "installment_payment"
$ws.Range("C4").Value = "Chords Music Academy (+917981585309)"
$ws.Range("D4").Value = "UTILITY"
$ws.Range("E4").Value = "APPROVED"
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = "Var1|Var2|Var3|Var4|Var5"
$ws.Range("H4").Value = "No media required"
$ws.Range("I4").Value = $apiUrl
$ws.Range("J4").Value = $msgContent

$ws.Rows.Item(4).AutoFit()
